# Swap the order of the two comma-separated entries in the "Recorded By"
# (column G) cells for a specific set of rows on the "Session Analysis
# Results" sheet, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row numbers (in column G) whose value needs to become the swapped order,
# together with the exact resulting text.
$updates = @{
    3   = "System, dnasr281@gmail.com"
    4   = "System, backup@backdoor.com"
    6   = "System, dnasr281@gmail.com"
    7   = "admin@admin.com, System"
    10  = "System, dnasr281@gmail.com"
    11  = "System, dnasr281@gmail.com"
    12  = "System, dnasr281@gmail.com"
    13  = "System, dnasr281@gmail.com"
    14  = "System, dnasr281@gmail.com"
    15  = "System, dnasr281@gmail.com"
    17  = "System, dnasr281@gmail.com"
    18  = "System, dnasr281@gmail.com"
    19  = "System, dnasr281@gmail.com"
    20  = "System, dnasr281@gmail.com"
    21  = "System, dnasr281@gmail.com"
    22  = "System, dnasr281@gmail.com"
    24  = "System, dnasr281@gmail.com"
    30  = "System, dnasr281@gmail.com"
    31  = "System, backup@backdoor.com"
    33  = "System, dnasr281@gmail.com"
    34  = "admin@admin.com, System"
    37  = "System, dnasr281@gmail.com"
    38  = "System, dnasr281@gmail.com"
    39  = "System, dnasr281@gmail.com"
    40  = "System, dnasr281@gmail.com"
    41  = "System, dnasr281@gmail.com"
    42  = "System, dnasr281@gmail.com"
    44  = "System, dnasr281@gmail.com"
    45  = "System, dnasr281@gmail.com"
    46  = "System, dnasr281@gmail.com"
    47  = "System, dnasr281@gmail.com"
    48  = "System, dnasr281@gmail.com"
    49  = "System, dnasr281@gmail.com"
    51  = "System, dnasr281@gmail.com"
    57  = "System, dnasr281@gmail.com"
    58  = "System, backup@backdoor.com"
    60  = "System, dnasr281@gmail.com"
    61  = "admin@admin.com, System"
    64  = "System, dnasr281@gmail.com"
    65  = "System, dnasr281@gmail.com"
    66  = "System, dnasr281@gmail.com"
    67  = "System, dnasr281@gmail.com"
    68  = "System, dnasr281@gmail.com"
    69  = "System, dnasr281@gmail.com"
    71  = "System, dnasr281@gmail.com"
    72  = "System, dnasr281@gmail.com"
    73  = "System, dnasr281@gmail.com"
    74  = "System, dnasr281@gmail.com"
    75  = "System, dnasr281@gmail.com"
    76  = "System, dnasr281@gmail.com"
    78  = "System, dnasr281@gmail.com"
    86  = "System, dnasr281@gmail.com"
    87  = "System, dnasr281@gmail.com"
    88  = "System, dnasr281@gmail.com"
    89  = "System, dnasr281@gmail.com"
    90  = "admin@admin.com, dnasr281@gmail.com"
    93  = "System, dnasr281@gmail.com"
    95  = "System, dnasr281@gmail.com"
    96  = "System, dnasr281@gmail.com"
    97  = "System, dnasr281@gmail.com"
    99  = "System, dnasr281@gmail.com"
    102 = "System, dnasr281@gmail.com"
    104 = "System, dnasr281@gmail.com"
    112 = "System, dnasr281@gmail.com"
    113 = "System, dnasr281@gmail.com"
    114 = "System, dnasr281@gmail.com"
    115 = "System, dnasr281@gmail.com"
    116 = "admin@admin.com, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    121 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    123 = "System, dnasr281@gmail.com"
    125 = "System, dnasr281@gmail.com"
    128 = "System, dnasr281@gmail.com"
    130 = "System, dnasr281@gmail.com"
    138 = "System, dnasr281@gmail.com"
    139 = "System, dnasr281@gmail.com"
    140 = "System, dnasr281@gmail.com"
    141 = "System, dnasr281@gmail.com"
    142 = "admin@admin.com, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
    147 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
    149 = "System, dnasr281@gmail.com"
    151 = "System, dnasr281@gmail.com"
    154 = "System, dnasr281@gmail.com"
    156 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
